# Apply the "Mantenimiento" budget-sheet addition described in the commit:
#   "Se añaden costes de mantenimiento #1"
#
# 1. Rename the original (only) sheet "Hoja1" -> "Desarrollo Proyecto"
# 2. Add a new sheet "Mantenimiento" right after it, with the same
#    header/row layout as the first sheet's "Personal" block, containing a
#    single maintenance-cost line item, a subtotal, and TOTAL (sin/con IVA)
#    rows.
# 3. Leave "Mantenimiento" as the active/selected sheet (tab shown in front).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Desarrollo Proyecto"

# New sheet, inserted immediately after "Desarrollo Proyecto"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Mantenimiento"

# ---- column widths (match "Desarrollo Proyecto" visual proportions) ------
$ws2.Columns.Item(1).ColumnWidth = 24.736979166666668
$ws2.Columns.Item(2).ColumnWidth = 53.877604166666664
$ws2.Columns.Item(3).ColumnWidth = 8.592447916666666
$ws2.Columns.Item(4).ColumnWidth = 24.307291666666668
$ws2.Columns.Item(5).ColumnWidth = 12.451822916666666

# ---- row heights -----------------------------------------------------
$ws2.Rows.Item(1).RowHeight = 15.75
$ws2.Rows.Item(5).RowHeight = 18.75
$ws2.Rows.Item(6).RowHeight = 18.75

# ---- formatting: clone styles from the matching cells on sheet 1 --------
# Header row (Tarea/Actividad, Elemento, Unidades, Coste por unidad, Coste Total)
$ws1.Range("A6:E6").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)  # xlPasteFormats

# Category row ("Personal"-style banner row)
$ws1.Range("A7:E7").Copy()
$ws2.Range("A2:E2").PasteSpecial(-4122)

# Item row (single line item)
$ws1.Range("A8:E8").Copy()
$ws2.Range("A3:E3").PasteSpecial(-4122)

# Subtotal row
$ws1.Range("A10:E10").Copy()
$ws2.Range("A4:E4").PasteSpecial(-4122)

# TOTAL (sin IVA) / TOTAL (con IVA 21%) rows
$ws1.Range("D21:E21").Copy()
$ws2.Range("D5:E5").PasteSpecial(-4122)
$ws1.Range("D22:E22").Copy()
$ws2.Range("D6:E6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- content -----------------------------------------------------------
$ws2.Range("A1").Value2 = "Tarea / Actividad"
$ws2.Range("B1").Value2 = "Elemento"
$ws2.Range("C1").Value2 = "Unidades"
$ws2.Range("D1").Value2 = "Coste por unidad"
$ws2.Range("E1").Value2 = "Coste Total"

$ws2.Range("A2").Value2 = "Mantenimiento"

$ws2.Range("B3").Value2 = "Disponivilidad y mantenimiento de los servidores (por mes)"
$ws2.Range("C3").Value2 = 1
$ws2.Range("D3").Value2 = 2500
$ws2.Range("E3").Formula = "=D3*C3"

$ws2.Range("D4").Value2 = "Subtotal 1"
$ws2.Range("E4").Formula = "=E3"

$ws2.Range("D5").Value2 = "TOTAL (sin IVA)"
$ws2.Range("E5").Formula = "=E4"

$ws2.Range("D6").Value2 = "TOTAL (con IVA 21%)"
$ws2.Range("E6").Formula = "=E5*1.21"

# ---- selection / active sheet -------------------------------------------
$ws1.Range("B22").Select()
$ws2.Range("B8").Select()
$ws2.Activate()
